$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp label (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 17:52"

# --- Country label swaps caused by re-ranking ---
# Japon overtook Corea del Sur (rows 26/27)
$ws.Range("A26").Value = "Japon"
$ws.Range("A27").Value = "Corea del Sur"

# Principado de Andorra overtook Crucero (rows 91/92)
$ws.Range("A91").Value = "Principado de Andorra"
$ws.Range("A92").Value = "Crucero"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 741230
$ws.Range("C4").Value = 2438
$ws.Range("D4").Value = 68610
$ws.Range("E4").Value = 633517
$ws.Range("G4").Value = 89
$ws.Range("H4").Value = 39103

# Row 16: Paises Bajos
$ws.Range("B16").Value = 33951
$ws.Range("C16").Value = 568
$ws.Range("E16").Value = 21235
$ws.Range("G16").Value = 39
$ws.Range("H16").Value = 1509

# Row 20: Peru
$ws.Range("B20").Value = 16960
$ws.Range("C20").Value = 595
$ws.Range("D20").Value = 2745
$ws.Range("E20").Value = 13665
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 550

# Row 22
$ws.Range("B22").Value = 14749
$ws.Range("C22").Value = 78
$ws.Range("E22").Value = 3796
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = 452

# Row 26: Japon (new data, moved up in ranking)
$ws.Range("B26").Value = 10797
$ws.Range("C26").Value = 501
$ws.Range("D26").Value = 1159
$ws.Range("E26").Value = 9416
$ws.Range("F26").Value = 217
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 222

# Row 27: Corea del Sur (unchanged data, moved down in ranking)
$ws.Range("B27").Value = 10661
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 8042
$ws.Range("E27").Value = 2385
$ws.Range("F27").Value = 55
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 234

# Row 28: Chile
$ws.Range("B28").Value = 10088
$ws.Range("C28").Value = 358
$ws.Range("D28").Value = 4338
$ws.Range("E28").Value = 5617
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 133

# Row 30: Ecuador
$ws.Range("B30").Value = 9287
$ws.Range("C30").Value = 545
$ws.Range("E30").Value = 7887
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = 360

# Row 40
$ws.Range("D40").Value = 768
$ws.Range("E40").Value = 5809
$ws.Range("F40").Value = 22

# Row 88
$ws.Range("D88").Value = 81
$ws.Range("E88").Value = 674
$ws.Range("F88").Value = 15

# Row 91: Principado de Andorra (new data, moved up in ranking)
$ws.Range("B91").Value = 713
$ws.Range("C91").Value = 9
$ws.Range("D91").Value = 235
$ws.Range("E91").Value = 442
$ws.Range("F91").Value = 17
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 36

# Row 92: Crucero (unchanged data, moved down in ranking)
$ws.Range("B92").Value = 712
$ws.Range("D92").Value = 644
$ws.Range("E92").Value = 55
$ws.Range("F92").Value = 7
$ws.Range("H92").Value = 13

# Row 102
$ws.Range("E102").Value = 209
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 10
